$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns ("sum_SASA", "max_SASA") before the current
# "flexibility" column (D), shifting D:J to F:L.
$ws.Range("D1:E1").EntireColumn.Insert()

# Drop the now-stale data rows (former rows 3-7); only one data row remains.
$ws.Range("A3:A7").EntireRow.Delete()

# Headers for the two newly inserted columns.
$ws.Range("D1").Value = "sum_SASA"
$ws.Range("E1").Value = "max_SASA"

# Replace the remaining data row with the updated values.
$ws.Range("A2").Value = "GlcNAc(b1-4)GlcNAc(b1-4)GlcNAc"
$ws.Range("B2").Value = 1.382957652431078
$ws.Range("C2").Value = 3.321409580704815
$ws.Range("D2").Value = 3.321409580704815
$ws.Range("E2").Value = 3.321409580704815
$ws.Range("F2").Value = 1.234
$ws.Range("G2").Value = 0.371
$ws.Range("H2").Value = 8.91
$ws.Range("I2").Value = "['4C1']"
$ws.Range("J2").Value = "['GlcNAc(b1-4)']"
$ws.Range("K2").Value = "['GlcNAc(b1-4)']"
$ws.Range("L2").Value = "N"
